$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $value)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "261.31"
Set-TextValue $ws "E2" "1.60%"
Set-TextValue $ws "D3" "27.41"
Set-TextValue $ws "E3" "1.30%"
Set-TextValue $ws "D4" "4.756"
Set-TextValue $ws "E4" "4.16%"
Set-TextValue $ws "D5" "0.06068"
Set-TextValue $ws "E5" "2.87%"
Set-TextValue $ws "D6" "6.638"
Set-TextValue $ws "E6" "0.13%"
Set-TextValue $ws "D7" "0.8622"
Set-TextValue $ws "E7" "0.87%"
Set-TextValue $ws "D8" "0.9205"
Set-TextValue $ws "E8" "-1.69%"
Set-TextValue $ws "D9" "0.1407"
Set-TextValue $ws "E9" "1.47%"
Set-TextValue $ws "D10" "0.05041"
Set-TextValue $ws "E10" "3.44%"
Set-TextValue $ws "D11" "0.07100"
Set-TextValue $ws "E11" "0.38%"
Set-TextValue $ws "D12" "0.03037"
Set-TextValue $ws "E12" "-0.91%"
Set-TextValue $ws "D13" "0.09093"
Set-TextValue $ws "E13" "-0.22%"
Set-TextValue $ws "D14" "0.001542"
Set-TextValue $ws "E14" "1.23%"
Set-TextValue $ws "D15" "0.0006087"
Set-TextValue $ws "E15" "0.94%"
Set-TextValue $ws "D16" "0.006164"
Set-TextValue $ws "E16" "2.26%"
Set-TextValue $ws "D17" "3.453"
Set-TextValue $ws "D18" "3.174"
Set-TextValue $ws "E18" "-0.26%"
Set-TextValue $ws "E20" "2.48%"
Set-TextValue $ws "D22" "4.101"
Set-TextValue $ws "E22" "4.85%"
Set-TextValue $ws "D23" "0.04238"
Set-TextValue $ws "E23" "-0.80%"
Set-TextValue $ws "D24" "0.001220"
Set-TextValue $ws "E24" "0.12%"
Set-TextValue $ws "E25" "-8.79%"
Set-TextValue $ws "E26" "0.03%"
Set-TextValue $ws "E27" "3.13%"
Set-TextValue $ws "D40" "0.03881"
Set-TextValue $ws "E40" "1.55%"
Set-TextValue $ws "E41" "1.22%"
Set-TextValue $ws "D42" "0.004131"
Set-TextValue $ws "E42" "-33.91%"
Set-TextValue $ws "D43" "0.01499"
Set-TextValue $ws "E43" "8.15%"
Set-TextValue $ws "D44" "0.002185"
Set-TextValue $ws "E44" "-0.67%"
Set-TextValue $ws "D45" "0.00005302"
Set-TextValue $ws "E45" "-1.28%"
Set-TextValue $ws "E46" "0.03%"
Set-TextValue $ws "E47" "-17.18%"
Set-TextValue $ws "E48" "-47.64%"
Set-TextValue $ws "E49" "0.03%"
Set-TextValue $ws "E50" "0.03%"
